# Update "想去人数" (want-to-go count) figures on the 展览 (sheet 1) and
# 全部类型 (sheet 4) worksheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F13").Value = 231
$ws1.Range("F14").Value = 579
$ws1.Range("F15").Value = 11577
$ws1.Range("F16").Value = 11764
$ws1.Range("F18").Value = 75

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F14").Value = 231
$ws4.Range("F15").Value = 579
$ws4.Range("F16").Value = 11577
$ws4.Range("F17").Value = 11764
$ws4.Range("F19").Value = 75
